# Fix the overvote delimiter used in the test CVR data: the delimiter
# character was changed from "/" to "|", so the shared-string values that
# encode overvotes (multiple candidates marked on one ranking) must be
# updated to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "D|A"
$ws.Range("C9").Value = "A|D"

# The workbook was also re-saved with a highlight fill behind the bordered
# data cells (white interior fill added to the existing bordered range),
# while the border itself keeps the same appearance.
$ws.Range("A1:E10").Interior.Color = 16777215
